$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 98 (shifts old rows 98..118 down to 100..120).
$ws.Range("A98:A99").EntireRow.Insert()

# New row 98: week of 2023-10-10 ("Primera" quality)
$ws.Range("A98").Value = 3
$ws.Range("B98").Value = "Femacal de La Calera"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 45209
$ws.Range("E98").Value = 5
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100108
$ws.Range("H98").Value = "Tropicales y subtropicales"
$ws.Range("I98").Value = 100108004
$ws.Range("J98").Value = "Papaya"
$ws.Range("K98").Value = "Cultivar IV Región"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 56
$ws.Range("N98").Value = 16000
$ws.Range("O98").Value = 16000
$ws.Range("P98").Value = 16000
$ws.Range("Q98").Value = "`$/bandeja 10 kilos"
$ws.Range("R98").Value = "Provincia del Elquí"
$ws.Range("S98").Value = 1600
$ws.Range("T98").Value = 10

# New row 99: week of 2023-10-10 ("Segunda" quality)
$ws.Range("A99").Value = 3
$ws.Range("B99").Value = "Femacal de La Calera"
$ws.Range("C99").Value = "Coquimbo"
$ws.Range("D99").Value = 45209
$ws.Range("E99").Value = 5
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100108
$ws.Range("H99").Value = "Tropicales y subtropicales"
$ws.Range("I99").Value = 100108004
$ws.Range("J99").Value = "Papaya"
$ws.Range("K99").Value = "Cultivar IV Región"
$ws.Range("L99").Value = "Segunda"
$ws.Range("M99").Value = 48
$ws.Range("N99").Value = 13000
$ws.Range("O99").Value = 13000
$ws.Range("P99").Value = 13000
$ws.Range("Q99").Value = "`$/bandeja 10 kilos"
$ws.Range("R99").Value = "Provincia del Elquí"
$ws.Range("S99").Value = 1300
$ws.Range("T99").Value = 10
